$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "trainingimages/10_tokiti"
$ws.Range("B2").Value = "pngimages/10_backpack.png"
$ws.Range("C2").Value = "trainingimages/19_papipi"
$ws.Range("D2").Value = "pngimages/19_burger.png"
$ws.Range("E2").Value = -0.5
$ws.Range("F2").Value = 0.5

# Row 3
$ws.Range("A3").Value = "trainingimages/09_tipata"
$ws.Range("B3").Value = "pngimages/09_plane.png"
$ws.Range("C3").Value = "trainingimages/22_kakoki"
$ws.Range("D3").Value = "pngimages/22_egg.png"
$ws.Range("E3").Value = -0.5
$ws.Range("F3").Value = 0.5

# Row 4
$ws.Range("A4").Value = "trainingimages/21_papika"
$ws.Range("B4").Value = "pngimages/21_cheese.png"
$ws.Range("C4").Value = "trainingimages/14_pokoto"
$ws.Range("D4").Value = "pngimages/14_coffee.png"
$ws.Range("E4").Value = 0.5
$ws.Range("F4").Value = -0.5
